# Update "想去人数" (F column) values across sheets to match the newly
# generated data snapshot (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6935
$ws1.Range("F4").Value = 50
$ws1.Range("F5").Value = 0
$ws1.Range("F6").Value = 157
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 200
$ws1.Range("F10").Value = 1290
$ws1.Range("F15").Value = 17
$ws1.Range("F16").Value = 406
$ws1.Range("F17").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("F23").Value = 0
$ws1.Range("F25").Value = 219

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 46

# --- Sheet "全部类型" (all types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 0
$ws4.Range("F5").Value = 453
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 0
$ws4.Range("F13").Value = 0
$ws4.Range("F14").Value = 145
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 0
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 27
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 5141
$ws4.Range("F23").Value = 0
$ws4.Range("F25").Value = 0
